# Apply the "456a3b4" data refresh to 广州-漫展信息.xlsx
#
# Sheet 1 "展览"   (exhibitions) : dimension A1:I29 -> A1:I30
# Sheet 2 "演出"   (performances): no row/dimension change
# Sheet 3 "本地生活"(local life) : untouched
# Sheet 4 "全部类型"(all types)  : dimension A1:I45 -> A1:I46

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a string value into a cell without letting the host
# auto-convert date-shaped text (e.g. "2024-07-06") into a real date
# serial number. Formatting the cell as Text first, then clearing the
# format again afterwards, keeps the stored cell a plain string while
# leaving no stray number-format behind.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# =======================================================================
# Sheet 1: 展览 (exhibitions)
# =======================================================================
$ws1 = $wb.Worksheets.Item(1)

# --- "想去人数" (F column) refresh on existing rows -------------------
$ws1.Cells.Item(5, 6).Value  = 1055
$ws1.Cells.Item(8, 6).Value  = 570
$ws1.Cells.Item(9, 6).Value  = 1490
$ws1.Cells.Item(11, 6).Value = 1392
$ws1.Cells.Item(12, 6).Value = 3039
$ws1.Cells.Item(13, 6).Value = 537
$ws1.Cells.Item(14, 6).Value = 1695
$ws1.Cells.Item(15, 6).Value = 1567
$ws1.Cells.Item(18, 6).Value = 1428
$ws1.Cells.Item(20, 6).Value = 68
$ws1.Cells.Item(22, 6).Value = 319
$ws1.Cells.Item(25, 6).Value = 3609
$ws1.Cells.Item(26, 6).Value = 714
$ws1.Cells.Item(28, 6).Value = 1593

# --- Insert a new event row before the old row 29 ----------------------
# Shift the existing row 29 (广州·火影only) down to row 30, preserving
# its formatting/types, then overwrite row 29 with the new event.
$ws1.Range("A29:I29").Copy($ws1.Range("A30:I30"))

# Fix up the shifted-down row: new index number + updated "想去人数"
$ws1.Cells.Item(30, 1).Value = 29
$ws1.Cells.Item(30, 6).Value = 57

# Fill the new row 29 with the "广州·重生之道only" event
Set-TextValue $ws1.Cells.Item(29, 2) "2024-07-06"
$ws1.Cells.Item(29, 3).Value = "广州·重生之道only"
$ws1.Cells.Item(29, 4).Value = "同泰路颐和山庄 颐和大酒店"
$ws1.Cells.Item(29, 5).Value = "2024.07.06 10:30-07.06 16:30"
$ws1.Cells.Item(29, 6).Value = 0
$ws1.Cells.Item(29, 7).Value = 75
$ws1.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84896"
$ws1.Cells.Item(29, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/aJpJGAEc1713699622756.png"

# =======================================================================
# Sheet 2: 演出 (performances) -- only two "想去人数" values change
# =======================================================================
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(9, 6).Value  = 35
$ws2.Cells.Item(13, 6).Value = 90

# =======================================================================
# Sheet 3: 本地生活 (local life) -- unchanged
# =======================================================================

# =======================================================================
# Sheet 4: 全部类型 (all types)
# =======================================================================
$ws4 = $wb.Worksheets.Item(4)

# --- "想去人数" (F column) refresh on existing rows -------------------
$ws4.Cells.Item(14, 6).Value = 35
$ws4.Cells.Item(16, 6).Value = 1055
$ws4.Cells.Item(19, 6).Value = 570
$ws4.Cells.Item(20, 6).Value = 1490
$ws4.Cells.Item(22, 6).Value = 1392
$ws4.Cells.Item(23, 6).Value = 3039
$ws4.Cells.Item(24, 6).Value = 537
$ws4.Cells.Item(25, 6).Value = 1695
$ws4.Cells.Item(26, 6).Value = 1567
$ws4.Cells.Item(29, 6).Value = 1428
$ws4.Cells.Item(31, 6).Value = 68
$ws4.Cells.Item(35, 6).Value = 319
$ws4.Cells.Item(38, 6).Value = 3609
$ws4.Cells.Item(39, 6).Value = 714
$ws4.Cells.Item(41, 6).Value = 1593
$ws4.Cells.Item(42, 6).Value = 90

# --- Insert a new event row before the old row 44 ----------------------
# Shift rows 45 -> 46 and 44 -> 45 (in that order) to make room, then
# overwrite row 44 in place with the new event.
$ws4.Range("A45:I45").Copy($ws4.Range("A46:I46"))
$ws4.Range("A44:I44").Copy($ws4.Range("A45:I45"))

# Fix up the shifted-down rows: new index numbers + updated "想去人数"
$ws4.Cells.Item(45, 1).Value = 44
$ws4.Cells.Item(45, 6).Value = 57
$ws4.Cells.Item(46, 1).Value = 45

# Fill row 44 with the "广州·重生之道only" event
Set-TextValue $ws4.Cells.Item(44, 2) "2024-07-06"
$ws4.Cells.Item(44, 3).Value = "广州·重生之道only"
$ws4.Cells.Item(44, 4).Value = "同泰路颐和山庄 颐和大酒店"
$ws4.Cells.Item(44, 5).Value = "2024.07.06 10:30-07.06 16:30"
$ws4.Cells.Item(44, 6).Value = 0
$ws4.Cells.Item(44, 7).Value = 75
$ws4.Cells.Item(44, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84896"
$ws4.Cells.Item(44, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/aJpJGAEc1713699622756.png"

Write-Host "Edit complete"
